$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B (47 -> 54 characters)
$ws.Columns.Item(2).ColumnWidth = 53.17

# Remove all existing hyperlinks (and their relationships) before rewriting the table
$ws.Cells.Hyperlinks.Delete()

# Rewrite the data rows (a new row was inserted at the top, another was inserted
# in the middle, and four new rows were appended at the bottom, with every
# timestamp bumped to the new scrape time)

# Row 2
$ws.Range("A2").Value = '2025-09-23 01:15:18'
$ws.Range("B2").Value = '初回 「AIで笑顔を検出し、2秒クリップを無劣化で自動切り出すWindowsツール開発(予算10万円)」'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5398662'
$ws.Range("G2").Value = 413
$ws.Range("H2").Value = '🔥AI,Ai ◆ツール,開発'

# Row 3
$ws.Range("A3").Value = '2025-09-23 01:15:18'
$ws.Range("B3").Value = '【急募】Pythonによるエキテンの自動スクレイピングツール開発依頼'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5398198'
$ws.Range("G3").Value = 343
$ws.Range("H3").Value = '🔥Python ◆ツール,開発'

# Row 4
$ws.Range("A4").Value = '2025-09-23 01:15:18'
$ws.Range("B4").Value = '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5217096'
$ws.Range("G4").Value = 243
$ws.Range("H4").Value = '🔥API ◆ツール'

# Row 5
$ws.Range("A5").Value = '2025-09-23 01:15:18'
$ws.Range("B5").Value = '【急募】出品・在庫管理ツール開発と保守対応者募集'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5398562'
$ws.Range("G5").Value = 163
$ws.Range("H5").Value = '◆ツール,開発 ◇管理'

# Row 6
$ws.Range("A6").Value = '2025-09-23 01:15:18'
$ws.Range("B6").Value = '【相談希望】在庫管理・出品補助ツールの開発に関するZoom面談依頼'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5398112'
$ws.Range("G6").Value = 158
$ws.Range("H6").Value = '◆ツール,開発 ◇管理'

# Row 7
$ws.Range("A7").Value = '2025-09-23 01:15:18'
$ws.Range("B7").Value = '【急募】Slack自動リアクションツール開発依頼'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5398193'
$ws.Range("G7").Value = 120
$ws.Range("H7").Value = '◆ツール,開発'

# Row 8
$ws.Range("A8").Value = '2025-09-23 01:15:18'
$ws.Range("B8").Value = '【急募】MT4特定口座の取引を子口座に反映するシステム開発'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5398203'
$ws.Range("G8").Value = 118
$ws.Range("H8").Value = '◆開発,システム開発'

# Row 9
$ws.Range("A9").Value = '2025-09-23 01:15:18'
$ws.Range("B9").Value = '【急募】自己分析アプリのバックエンド開発アドバイザリー募集'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5397930'
$ws.Range("G9").Value = 93
$ws.Range("H9").Value = '◆開発 ◇アプリ'

# Row 10
$ws.Range("A10").Value = '2025-09-23 01:15:18'
$ws.Range("B10").Value = '【GAS開発者募集】Amazon広告管理SaaSのMVP開発'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5397812'
$ws.Range("G10").Value = 88
$ws.Range("H10").Value = '◆開発 ◇管理'

# Row 11
$ws.Range("A11").Value = '2025-09-23 01:15:18'
$ws.Range("B11").Value = 'EC多プラットフォーム展開在庫・価格連携ツールの作成'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5398432'
$ws.Range("G11").Value = 80
$ws.Range("H11").Value = '◆ツール'

# Row 12
$ws.Range("A12").Value = '2025-09-23 01:15:18'
$ws.Range("B12").Value = '【急募】iOSアプリのAdMobメディエーション入札接続とeCPM改善'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5398081'
$ws.Range("G12").Value = 38
$ws.Range("H12").Value = '◇アプリ'

# Row 13
$ws.Range("A13").Value = '2025-09-23 01:15:18'
$ws.Range("B13").Value = '【急募】PHP・Lalavelでの既存プログラム改修依頼'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5396563'
$ws.Range("G13").Value = 33
$ws.Range("H13").Value = '○PHP'

# Row 14
$ws.Range("A14").Value = '2025-09-23 01:15:18'
$ws.Range("B14").Value = '初回 iOSとAndroidのアプリ 課金(サブスク)'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5398382'
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = '◇アプリ'

# Row 15
$ws.Range("A15").Value = '2025-09-23 01:15:18'
$ws.Range("B15").Value = '【急募】災害時に備えた「ピジョンクラウド」でのシステムづくり、運用サポートの依頼'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5398657'
$ws.Range("G15").Value = 33
$ws.Range("H15").ClearContents()

# Row 16
$ws.Range("A16").Value = '2025-09-23 01:15:18'
$ws.Range("B16").Value = '【Braze経験者募集】CRM/マーケティングオートメーション支援(中級者以上)'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5398071'
$ws.Range("G16").Value = 25
$ws.Range("H16").ClearContents()

# Row 17
$ws.Range("A17").Value = '2025-09-23 01:15:18'
$ws.Range("B17").Value = '【Braze経験者募集】CRM/マーケティングオートメーション支援(中級者以上)'
$ws.Range("C17").Value = 'システム開発'
$ws.Range("D17").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E17").Value = '期限情報なし'
$ws.Range("F17").Value = 'https://www.lancers.jp/work/detail/5398062'
$ws.Range("G17").Value = 25
$ws.Range("H17").ClearContents()

# Row 18
$ws.Range("A18").Value = '2025-09-23 01:15:18'
$ws.Range("B18").Value = '限定公開 限定公開の仕事'
$ws.Range("C18").Value = 'システム開発'
$ws.Range("D18").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E18").Value = '期限情報なし'
$ws.Range("F18").Value = 'https://www.lancers.jp/work/detail/5398293'
$ws.Range("G18").Value = 18
$ws.Range("H18").ClearContents()

# Row 19
$ws.Range("A19").Value = '2025-09-23 01:15:18'
$ws.Range("B19").Value = 'データセンター向けサーバー・ルーター設置作業'
$ws.Range("C19").Value = 'システム開発'
$ws.Range("D19").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E19").Value = '期限情報なし'
$ws.Range("F19").Value = 'https://www.lancers.jp/work/detail/5397887'
$ws.Range("G19").Value = 18
$ws.Range("H19").ClearContents()

# Row 20
$ws.Range("A20").Value = '2025-09-23 01:15:18'
$ws.Range("B20").Value = 'Excelでのデータシュミレーション'
$ws.Range("C20").Value = 'システム開発'
$ws.Range("D20").Value = '~ 5,000 円 / 固定'
$ws.Range("E20").Value = '期限情報なし'
$ws.Range("F20").Value = 'https://www.lancers.jp/work/detail/5398497'
$ws.Range("G20").Value = 10
$ws.Range("H20").ClearContents()

# Row 21
$ws.Range("A21").Value = '2025-09-23 01:15:18'
$ws.Range("B21").Value = 'Excelやスプレッドシートでのデータシュミレーション クエリ(query)や関数利用'
$ws.Range("C21").Value = 'システム開発'
$ws.Range("D21").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E21").Value = '期限情報なし'
$ws.Range("F21").Value = 'https://www.lancers.jp/work/detail/5397980'
$ws.Range("G21").Value = 10
$ws.Range("H21").ClearContents()

# Row 22
$ws.Range("A22").Value = '2025-09-23 01:15:18'
$ws.Range("B22").Value = '【中小企業支援】債務超過・赤字経営解消の診断依頼'
$ws.Range("C22").Value = 'システム開発'
$ws.Range("D22").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E22").Value = '期限情報なし'
$ws.Range("F22").Value = 'https://www.lancers.jp/work/detail/5397962'
$ws.Range("G22").Value = 10
$ws.Range("H22").ClearContents()

# Row 23
$ws.Range("A23").Value = '2025-09-23 01:15:18'
$ws.Range("B23").Value = 'Geminiで旅行のしおりのHTMLを生成するプロンプトの作成'
$ws.Range("C23").Value = 'システム開発'
$ws.Range("D23").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E23").Value = '期限情報なし'
$ws.Range("F23").Value = 'https://www.lancers.jp/work/detail/5397817'
$ws.Range("G23").Value = 10
$ws.Range("H23").ClearContents()

# Re-create the hyperlinks on column F, in row order, then restore the "Hyperlink"
# cell style (Hyperlinks.Add alone creates a fresh duplicate style entry)
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5398662')
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5398198')
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5217096')
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5398562')
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5398112')
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5398193')
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5398203')
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5397930')
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5397812')
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5398432')
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5398081')
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5396563')
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5398382')
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5398657')
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5398071')
$ws.Hyperlinks.Add($ws.Range("F17"), 'https://www.lancers.jp/work/detail/5398062')
$ws.Hyperlinks.Add($ws.Range("F18"), 'https://www.lancers.jp/work/detail/5398293')
$ws.Hyperlinks.Add($ws.Range("F19"), 'https://www.lancers.jp/work/detail/5397887')
$ws.Hyperlinks.Add($ws.Range("F20"), 'https://www.lancers.jp/work/detail/5398497')
$ws.Hyperlinks.Add($ws.Range("F21"), 'https://www.lancers.jp/work/detail/5397980')
$ws.Hyperlinks.Add($ws.Range("F22"), 'https://www.lancers.jp/work/detail/5397962')
$ws.Hyperlinks.Add($ws.Range("F23"), 'https://www.lancers.jp/work/detail/5397817')

$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("F9").Style = "Hyperlink"
$ws.Range("F10").Style = "Hyperlink"
$ws.Range("F11").Style = "Hyperlink"
$ws.Range("F12").Style = "Hyperlink"
$ws.Range("F13").Style = "Hyperlink"
$ws.Range("F14").Style = "Hyperlink"
$ws.Range("F15").Style = "Hyperlink"
$ws.Range("F16").Style = "Hyperlink"
$ws.Range("F17").Style = "Hyperlink"
$ws.Range("F18").Style = "Hyperlink"
$ws.Range("F19").Style = "Hyperlink"
$ws.Range("F20").Style = "Hyperlink"
$ws.Range("F21").Style = "Hyperlink"
$ws.Range("F22").Style = "Hyperlink"
$ws.Range("F23").Style = "Hyperlink"
